$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "CasesTab" query (cell B2) dropped its trailing `Cohort` column —
# the MATCH (co:cohort) clause and the `co` bindings stay, only the final
# RETURN line that projected `co.cohort_description AS Cohort` is removed,
# and the preceding line's trailing comma goes away since it is now last.
$casesQuery = @'
MATCH (s:study)<-[*]-(c:case)<--(demo:demographic)
MATCH (c)<--(diag:diagnosis)
MATCH (co:cohort)<-[*]-(c)
WHERE diag.disease_term IN ['T Cell Lymphoma']
WITH DISTINCT c, s, demo, diag, co
RETURN  coalesce(c.case_id, '') AS `Case ID` ,
        coalesce(s.clinical_study_designation, '') AS `Study Code` ,
        coalesce(s.clinical_study_type, '') AS  `Study Type`,
        coalesce(demo.breed, '') AS Breed ,
        coalesce(diag.disease_term,'') AS Diagnosis ,
        coalesce(diag.stage_of_disease, '') AS `Stage of Disease` ,
        coalesce(demo.patient_age_at_enrollment, '') AS Age ,
        coalesce(demo.sex, '') AS Sex ,
        coalesce(demo.neutered_indicator, '') AS `Neutered Status`,
        coalesce(demo.weight, '') AS `Weight (kg)`,
        coalesce(diag.best_response, '') AS `Response to Treatment`
'@

$ws.Cells.Item(2, 2).Value = $casesQuery.TrimEnd("`r", "`n")

# Selection / view moved from the bottom row (B4) back to the top data row (B2),
# and the zoom was reset from 55% to 100%.
[void]$ws.Range("B2").Select()
$excel.ActiveWindow.Zoom = 100
